$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '97.125.89'
$ws.Range('E2').Value = '  +2.28%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.587.91'
$ws.Range('E3').Value = '  +0.95%  '

# Row 4
$ws.Range('E4').Value = '  +0.02%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.95'
$ws.Range('E5').Value = '  +2.32%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '654.72'
$ws.Range('E6').Value = '  +0.97%  '

# Row 7
$ws.Range('E7').Value = '  +15.63%  '

# Row 8
$ws.Range('E8').Value = '  +6.15%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  -0.04%  '

# Row 10
$ws.Range('E10').Value = '  +4.99%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.581.53'
$ws.Range('E11').Value = '  +0.84%  '

# Row 12
$ws.Range('E12').Value = '  +5.12%  '

# Row 13
$ws.Range('E13').Value = '  +1.02%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.39'
$ws.Range('E14').Value = '  -1.13%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.255.75'
$ws.Range('E15').Value = '  -0.04%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '97.270.93'
$ws.Range('E16').Value = '  +2.44%  '

# Row 17
$ws.Range('E17').Value = '  +2.84%  '

# Row 18
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.592.42'
$ws.Range('E18').Value = '  +1.20%  '

# Row 19
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.13'
$ws.Range('E19').Value = '  +2.63%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.68'
$ws.Range('E20').Value = '  +0.98%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '18.01'
$ws.Range('E21').Value = '  +1.07%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.520'
$ws.Range('E22').Value = '  +9.35%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.49'
$ws.Range('E23').Value = '  +1.03%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '513.55'
$ws.Range('E24').Value = '  +1.72%  '

# Row 25
$ws.Range('E25').Value = '  +5.83%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.51'
$ws.Range('E26').Value = '  -3.75%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '100.83'
$ws.Range('E27').Value = '  +6.14%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.99'
$ws.Range('E28').Value = '  +4.56%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.784.32'
$ws.Range('E29').Value = '  +1.03%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.160'
$ws.Range('E30').Value = '  +13.43%  '

# Row 31
$ws.Range('E31').Value = '  -0.30%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.79'
$ws.Range('E32').Value = '  +3.73%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  +0.02%  '

# Row 34
$ws.Range('E34').Value = '  +3.54%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.997'
$ws.Range('E35').Value = '  -0.13%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '31.65'
$ws.Range('E36').Value = '  +0.04%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '619.55'
$ws.Range('E37').Value = '  +5.91%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.73'
$ws.Range('E38').Value = '  +3.17%  '

# Row 39
$ws.Range('E39').Value = '  +1.35%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.64'
$ws.Range('E40').Value = '  +2.27%  '

# Row 41
$ws.Range('E41').Value = '  +11.76%  '

# Row 42
$ws.Range('E42').Value = '  +2.98%  '

# Row 43
$ws.Range('E43').Value = '  -0.06%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.921'
$ws.Range('E44').Value = '  +2.46%  '

# Row 45
$ws.Range('E45').Value = '  +5.54%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0439'
$ws.Range('E46').Value = '  +6.80%  '

# Row 47
$ws.Range('E47').Value = '  +0.37%  '

# Row 48
$ws.Range('E48').Value = '  +1.12%  '

# Row 49
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.402'
$ws.Range('E49').Value = '  +33.01%  '

# Row 50
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '33.10'
$ws.Range('E50').Value = '  -0.67%  '

# Row 51
$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.46'
$ws.Range('E51').Value = '  +5.13%  '
